$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Theme Success" to "Theme"
$ws.Name = "Theme"

# Add a new "Created" label in column A for rows 2-4
$ws.Range("A2").Value = "Created"
$ws.Range("A3").Value = "Created"
$ws.Range("A4").Value = "Created"

# Match the resulting selection/cursor position left behind in the saved file
$ws.Range("I16").Select()
